# Common: Added boring stuff for liquid
#
# Adds six new vendor names to the "vendors" sheet (column A, below the
# header in row 1) and re-sorts the list alphabetically, exactly as a
# human editor would do in Excel: type the new names under the existing
# list, then re-apply Data > Sort (A to Z) on the whole range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")

# Existing data occupies A2:A75 (74 vendors). Append the six new vendors
# right after the last existing row; the subsequent Sort call will move
# them into their correct alphabetical position.
$ws.Range("A76").Value = "Tobacco Bastards"
$ws.Range("A77").Value = "Rocket Girl"
$ws.Range("A78").Value = "Electra"
$ws.Range("A79").Value = "Five Pawns"
$ws.Range("A80").Value = "MaZa"
$ws.Range("A81").Value = "Ripe Vapes"

# Re-sort the full vendor list (A2:A81) alphabetically, same as the
# original "sortState" already stored on the sheet.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A81"))
$ws.Sort.SetRange($ws.Range("A2:A81"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Update the on-screen selection to match where the editor ended up.
$ws.Range("A56").Select() | Out-Null
